$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("YDS")
$ws.Range("B2").Value = '0 11 9 1 3 3 38 6 11 3 3 3 2 10 2 7 3 2 1 7 -1 0 0 6 0 1 1 -4 -1 3 3 0 1 2 5 2 7 4 5 2 2 10 2 3 2 7 4 21 2 0 1 3 9 3 0 8 -1 8 1 1 3 4 -3 0 9 2 2 0 1 2 2 16 8 4 0 7 9 5 1 1 3 -1 5 6 9 2 7 4 2 8 4 5 5 -1 5 15 7 -1 1 6 7 2 4 3 3 6 1 7 7 14 0 2 1 -3 5 11 9 9 7 7 14 4 -2 1 7 7 12 1 13 -2 15 1 5 6 7 4 5 3 20 14 14 12 6 2 3 6 5 3 -1 2 3 -3 2 0 4 1 5 4 7 1 1 27 1 -1 3 0 4 2 3 11 5 27 1 5 39 -1 12 -2 0 1 9 -5 -2 11 14 5 2 -1 2 3 4 7 5 1 9 4 0 2 0 8 6 5 -1 9 -1 1 6 2 13 2 8 17 4 -2 4 4 4 -6 3 -3 4 8 4 2 0 -2 2 1 9 2 0 -4 8 4 4 1 2 2 2 4 0 12 3 2 9 0 1 3 11 -1 25 3 5 -1 1 0 4 4 3 11 21 -2 12 0 5 2 4 4 -1 3 2 3 7 4 3 -3 6 15 20 10 1 3 4 9 3 4 10 4 4 2 5 49 4 3 2 4 12 0 1 3 -6 5 1 0 6 7 9 4 3 1 1 -1 4 5 0 6 2 7 5 2 9 2 4 0 2 2 5 4 0 2 6 3 1'
$ws.Range("C2").Value = '10 2 4 9 20 7 7 1 1 3 16 15 -1 -4 3 -4 1 3 6 1 10 1 2 2 0 5 10 5 2 6 6 12 -2 9 -2 -3 1 7 9 7 6 6 13 3 3 0 8 9 27 1 2 5 -4 1 3 2 4 2 10 4 2 6 1 6 3 7 2 6 2 12 4 2 12 4 2 2 4 2 -2 -1 5 4 3 5 2 12 -1 6 8 6 16 1 6 4 3 -2 4 14 0 1 5 2 4 5 26 1 3 9 1 5 0 12 1 0 -2 11 4 1 2 10 1 2 8 4 2 6 2 0 33 -2 4 16 1 1 5 10 -2 0 1 9 8 3 3 0 5 0 8 5 17 -1 0 5 0 5 3 -1 -1 4 22 9 -2 1 3 6 9 6 3 2 3 4 16 3 9 3 16 2 0 15 -3 -4 4 -4 4 6 3 0 14 8 4 0 0 3 8 1 4 1 4 2 1 1 2 3 3 1 -4 1 5 8 2 35 21 5 1 7 6 2 4 1 3 5 2 3 -3 2 6 -2 10 6 1 3 11 4 3 6 7 -1 1 7 10 -3 2 7 2 2 4 4 4 4 4 6 1 4 1 3 1 1 0 4 7 4 30 0 6 -4 7 1 3 1 2 -5 0 11 12 0 7 4 2 1 4 2 3 3 3 3 2 0 4 1 10 4 2 0 -1'
$ws.Range("B3").Value = '23 7 5 29 8 7 23 7 14 18 35 8 79 16 2 15 12 11 10 1 8 7 6 11 17 5 8 5 6 2 40 11 9 4 6 6 6 0 15 8 8 7 6 4 16 8 16 10 11 5 7 19 1 0 16 8 0 16 12 9 3 39 12 12 12 9 8 21 9 18 21 7 15 14 26 7 10 8 8 14 10 12 12 7 8 16 14 26 6 1 18 6 12 4 26 5 13 26 13 9 9 0 28 10 10 16 6 1 29 27 14 3 9 5 5 21 -4 19 11 6 23 16 6 15 50 3 83 15 7 11 19 23 4 19 18 16 7 11 14 22 11 10 16 6 5 19 29 1 11 32 7 10 15 6 5 12 5 3 7 5 -1 18 6 7 8 10 20 21 6 3 19 7 11 7 9 7 8 13 34 7 6 9 9 23 15 7 13 4 11 21 37 3 2 11 2 9 24 3 13 13 30 12 8 11 5 26 24 8 28 12 4 7 9 8 9 48 3 21 13 29 11 18 14 10 7 16'
$ws.Range("C3").Value = '22 6 -8 12 12 1 3 6 2 6 7 10 -1 7 11 6 43 3 23 5 5 3 10 5 6 5 3 8 7 9 11 13 3 20 2 20 25 4 11 4 9 1 5 26 91 6 10 9 14 42 11 -4 2 12 2 47 14 1 0 10 3 4 15 19 6 1 9 6 7 12 25 17 -2 5 3 2 28 28 12 1 19 9 2 13 16 6 2 5 7 8 13 10 9 23 33 3 6 2 9 14 8 5 8 13 7 9 2 30 9 11 -1 57 11 -2 2 -1 -1 5 5 -2 2 14 13 9 5 6 28 14 6 11 11 16 6 8 10 4 2 -6 19 22 6 2 9 10 20 5 10 0 9 9 8 0 21 17 33 50 17 6 -2 25 45 4 4 1 5 5 9 6 0 5 2 6 8 16 29 13 10 8 33 5 -2 5 2 6 -3 12 16 6 11 5 9 5 2 25 9 6 6 18 16 6 8 15 7 3 7 11 13 22 8 9 2 11 29 7 4 30 15 20 7 6 14 15 24 16 5 17 7 5 23 2 -1 4 1 0 8 10 7 12 -1 33 5 9 2 4 11 7 6 7 16 17 5 7 2 12 11 5 2 27 3 0'

$ws = $wb.Worksheets.Item("OFF")
$ws.Range("C2").Value = 170
$ws.Range("F2").Value = 46
$ws.Range("G2").Value = 48
$ws.Range("J2").Value = 21
$ws.Range("N2").Value = 6
$ws.Range("O2").Value = 16
$ws.Range("P2").Value = 9
$ws.Range("B3").Value = 11
$ws.Range("C3").Value = 104
$ws.Range("E3").Value = 26
$ws.Range("F3").Value = 65
$ws.Range("G3").Value = 19
$ws.Range("H3").Value = 22
$ws.Range("I3").Value = 39
$ws.Range("J3").Value = 28
$ws.Range("L3").Value = 163
$ws.Range("M3").Value = 108
$ws.Range("Q3").Value = 351

$ws = $wb.Worksheets.Item("DEF")
$ws.Range("C2").Value = 140
$ws.Range("D2").Value = 9
$ws.Range("E2").Value = 12
$ws.Range("F2").Value = 45
$ws.Range("G2").Value = 38
$ws.Range("H2").Value = 4
$ws.Range("J2").Value = 18
$ws.Range("N2").Value = 18
$ws.Range("O2").Value = 19
$ws.Range("P2").Value = 8
$ws.Range("B3").Value = 7
$ws.Range("C3").Value = 131
$ws.Range("E3").Value = 26
$ws.Range("F3").Value = 68
$ws.Range("G3").Value = 24
$ws.Range("H3").Value = 22
$ws.Range("I3").Value = 42
$ws.Range("J3").Value = 42
$ws.Range("L3").Value = 201
$ws.Range("M3").Value = 137
$ws.Range("Q3").Value = 375

$ws = $wb.Worksheets.Item("ST")
$ws.Range("B2").Value = 62
$ws.Range("D2").Value = 36
$ws.Range("F2").Value = 574
$ws.Range("G2").Value = 559
$ws.Range("N2").Value = 49
$ws.Range("O2").Value = 34
$ws.Range("B4").Value = '67 64 56 65 63 62 64 52 70 54 68 73 63 61 60 66 63 58 66 57 61 60 67 70 64 51 55 65 64 64 63 66 61 56 69'
$ws.Range("D4").Value = '0 0 11 2 0 0 0 0 0 0 6 0 0 0 11 0 12 0 5 0 0 0 0 0 2 0 0 0 0 0 0 0 0 0 0 11'
$ws.Range("B5").Value = '27 19 0 20 21 33 18 14 31 10 34 24 30 27 21 29 24 22 26 17 17 16 23 28 21 18 30 21 99 22 23 25 33 26 19'
$ws.Range("D5").Value = '0 0'
$ws.Range("D3").Value = '43 33 48 52 49 35 42 59 56 41 57 45 54 37 45 33 58 39 42 65 48 38 34 44 57 28 41 46 40 45 65 51 67 60 37 40'
$ws.Range("B6").Value = '16 19 16 19 11 68 16 20 26 7 21 10 22 27 27 33 29 22 23 37 0 21 12 12 23'

$ws = $wb.Worksheets.Item("TURNS")
$ws.Range("B3").Value = 3
$ws.Range("C3").Value = 3
$ws.Range("E3").Value = 14

$ws = $wb.Worksheets.Item("PEN")
$ws.Range("A2").Value = 'False start'
$ws.Range("B2").Value = 13
$ws.Range("C2").Value = 'Offside'
$ws.Range("D2").Value = 4
$ws.Range("A3").Value = 'Holding'
$ws.Range("B3").Value = 13
$ws.Range("C3").Value = 'Holding'
$ws.Range("D3").Value = 2
$ws.Range("A4").Value = 'Pass interference'
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 'Pass interference'
$ws.Range("D4").Value = 16
$ws.Range("A5").Value = 'Intentional grounding'
$ws.Range("B5").Value = 1

# Selection state to match target workbook view
$wsTurns = $wb.Worksheets.Item("TURNS")
$wsTurns.Range("C1").Select()
$wsSt = $wb.Worksheets.Item("ST")
$wsSt.Activate()
$wsSt.Range("D6").Select()